# Flask-calendar template: switch "Recuperativo"/"Especial" columns for
# "Extraordinario 1"/"Extraordinario 2", clear the per-subject names from
# the body rows (the generator now fills these in, along with the exam
# venue/time that used to be missing), add a new blank info row below the
# table, and nudge the column widths to fit the new headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename last two result columns -------------------------
# (the backing table's ListColumns/sharedStrings follow the header cells
# automatically, so this alone renames the Excel Table too)
$ws.Range("F6").Value = "Extraordinario 1"
$ws.Range("G6").Value = "Extraordinario 2"

# --- Body rows: the per-subject names are no longer hard-coded ----------
$ws.Range("A9:A15").ClearContents()

# --- New row under the table for venue/time info -------------------------
$ws.Range("A15:G15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false | Out-Null
$ws.Rows.Item(16).RowHeight = 50.25

# --- Column widths: widen to fit the new headers/content -----------------
$ws.Columns.Item(2).ColumnWidth = 27.67
$ws.Columns.Item(3).ColumnWidth = 30.5
$ws.Columns.Item(4).ColumnWidth = 31.83
$ws.Columns.Item(5).ColumnWidth = 31.83
$ws.Columns.Item(6).ColumnWidth = 30.67
$ws.Columns.Item(7).ColumnWidth = 46.67

# --- Selection as left by the editing session -----------------------------
$ws.Range("F9").Select() | Out-Null
